$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values for the columns that vary by row.
# Rows 2-9 hold the data; columns D, K, L, M, N, O, P, R, S are the ones
# that differ from row to row.
$cols = @("D","K","L","M","N","O","P","R","S")
$rows = 2..9

$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# The rows got reshuffled: each new row takes on the values that used to
# belong to a different row (a permutation of the 8 data rows).
$mapping = @{
    2 = 6
    3 = 7
    4 = 8
    5 = 2
    6 = 9
    7 = 5
    8 = 3
    9 = 4
}

foreach ($newRow in $rows) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
